$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G5").Value = "(words, that), (should, fail)"
$ws.Range("G6").Select()
